$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$ws.Range("D2").Value = "23.190.61"
$ws.Range("E2").Value = "  +0.37%  "
$ws.Range("D3").Value = "1.601.65"
$ws.Range("E3").Value = "  -0.03%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("E5").Value = "  -0.01%  "
Set-TextValue $ws.Range("D6") "303.36"
$ws.Range("E6").Value = "  +0.60%  "
Set-TextValue $ws.Range("D7") "0.3782"
$ws.Range("E7").Value = "  -0.06%  "
Set-TextValue $ws.Range("D8") "52.00"
$ws.Range("E8").Value = "  +4.42%  "
Set-TextValue $ws.Range("D9") "0.3615"
$ws.Range("E9").Value = "  -0.94%  "
Set-TextValue $ws.Range("D10") "1.269"
$ws.Range("E10").Value = "  +0.28%  "
$ws.Range("E11").Value = "  +0.00%  "
Set-TextValue $ws.Range("D12") "0.08115"
$ws.Range("E12").Value = "  -0.59%  "
Set-TextValue $ws.Range("D13") "22.74"
$ws.Range("E13").Value = "  -1.61%  "
Set-TextValue $ws.Range("D14") "6.580"
$ws.Range("E14").Value = "  -0.21%  "
Set-TextValue $ws.Range("D15") "7.411"
$ws.Range("E15").Value = "  +0.19%  "
$ws.Range("E16").Value = "  -1.42%  "
$ws.Range("D17").Value = "1.602.86"
$ws.Range("E17").Value = "  +0.22%  "
Set-TextValue $ws.Range("D18") "94.16"
$ws.Range("E18").Value = "  +2.47%  "
Set-TextValue $ws.Range("D19") "0.06880"
$ws.Range("E19").Value = "  +0.43%  "
Set-TextValue $ws.Range("D20") "18.09"
$ws.Range("E20").Value = "  -1.00%  "
Set-TextValue $ws.Range("D21") "6.546"
$ws.Range("E21").Value = "  -0.63%  "
Set-TextValue $ws.Range("D22") "1.000"
$ws.Range("E22").Value = "  -0.01%  "
Set-TextValue $ws.Range("D23") "12.98"
$ws.Range("E23").Value = "  -0.50%  "
$ws.Range("D24").Value = "23.182.72"
$ws.Range("E24").Value = "  +0.33%  "
$ws.Range("E25").Value = "  +2.54%  "
Set-TextValue $ws.Range("D26") "2.981"
$ws.Range("E26").Value = "  +9.46%  "
$ws.Range("E27").Value = "  +0.33%  "
Set-TextValue $ws.Range("D28") "149.44"
$ws.Range("E28").Value = "  -0.36%  "
Set-TextValue $ws.Range("D29") "5.242"
$ws.Range("E29").Value = "  -0.59%  "
Set-TextValue $ws.Range("D30") "133.96"
$ws.Range("E30").Value = "  +0.94%  "
Set-TextValue $ws.Range("D31") "2.387"
$ws.Range("E31").Value = "  -0.20%  "
Set-TextValue $ws.Range("D32") "6.760"
$ws.Range("E32").Value = "  -1.52%  "
$ws.Range("D33").Value = "1.779.49"
$ws.Range("E33").Value = "  +0.20%  "
Set-TextValue $ws.Range("D34") "0.9688"
$ws.Range("E34").Value = "  +0.89%  "
Set-TextValue $ws.Range("D35") "0.07506"
$ws.Range("E35").Value = "  -2.55%  "
Set-TextValue $ws.Range("D36") "10.30"
$ws.Range("E36").Value = "  +2.23%  "
Set-TextValue $ws.Range("D37") "0.02714"
$ws.Range("E37").Value = "  -0.38%  "
Set-TextValue $ws.Range("D38") "0.2502"
$ws.Range("E38").Value = "  -2.03%  "
Set-TextValue $ws.Range("D39") "0.08804"
$ws.Range("E39").Value = "  -0.93%  "
Set-TextValue $ws.Range("D40") "6.088"
$ws.Range("E40").Value = "  -3.01%  "
Set-TextValue $ws.Range("D41") "0.7113"
$ws.Range("E41").Value = "  +0.53%  "
Set-TextValue $ws.Range("D42") "1.359"
$ws.Range("E42").Value = "  -0.79%  "
Set-TextValue $ws.Range("D43") "12.49"
$ws.Range("E43").Value = "  -0.86%  "
Set-TextValue $ws.Range("D44") "15.61"
$ws.Range("E44").Value = "  +1.97%  "
Set-TextValue $ws.Range("D45") "0.6526"
$ws.Range("E45").Value = "  -1.27%  "
Set-TextValue $ws.Range("D46") "2.313"
$ws.Range("E46").Value = "  -0.07%  "
$ws.Range("E47").Value = "  +0.59%  "
Set-TextValue $ws.Range("D48") "132.10"
$ws.Range("E48").Value = "  -0.22%  "
$ws.Range("E49").Value = "  +0.35%  "
$ws.Range("E50").Value = "  -2.59%  "
Set-TextValue $ws.Range("D51") "1.213"
$ws.Range("E51").Value = "  +1.37%  "
